$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Establish cell formatting (borders/number formats) for the new rows by
#    copying from existing rows that already carry the correct alternating
#    style pattern used throughout this "form responses" sheet.
# ---------------------------------------------------------------------------

# Row 295 currently holds the special "last row" formatting (including the
# trailing, border-only M/N cells). Move that formatting down to the new
# last row (312) before anything else, so row 295 can go back to being a
# normal interior row.
$ws.Range("A295:N295").Copy()
$ws.Range("A312:N312").PasteSpecial(-4122)

# Rows 296-311 use the regular alternating style pairs found elsewhere in
# the table (A-L columns only - M/N are filled in separately below only
# where the new data actually needs them).
$ws.Range("A292:L293").Copy()
$ws.Range("A296:L297").PasteSpecial(-4122)
$ws.Range("A298:L299").PasteSpecial(-4122)
$ws.Range("A300:L301").PasteSpecial(-4122)
$ws.Range("A302:L303").PasteSpecial(-4122)
$ws.Range("A304:L305").PasteSpecial(-4122)
$ws.Range("A306:L307").PasteSpecial(-4122)
$ws.Range("A308:L309").PasteSpecial(-4122)
$ws.Range("A310:L311").PasteSpecial(-4122)

# A few of the new rows have non-blank M or N cells; pull the matching
# single-column style from existing rows that already use it.
$ws.Range("M291").Copy()
$ws.Range("M299").PasteSpecial(-4122)

$ws.Range("M280").Copy()
$ws.Range("M302").PasteSpecial(-4122)
$ws.Range("M306").PasteSpecial(-4122)

$ws.Range("N289").Copy()
$ws.Range("N301").PasteSpecial(-4122)
$ws.Range("N307").PasteSpecial(-4122)
$ws.Range("N309").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Populate the new rows (296-312) with the survey response data.
# ---------------------------------------------------------------------------
# Row 296
$ws.Cells.Item(296,1).Value = 45599.6874025
$ws.Cells.Item(296,2).Value = "seungye04@naver.com"
$ws.Cells.Item(296,3).Value = "언론방송융합미디어"
$ws.Cells.Item(296,4).Value = 20233846
$ws.Cells.Item(296,5).Value = "정승예"
$ws.Cells.Item(296,6).Value = "대한민국"
$ws.Cells.Item(296,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(296,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(296,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(296,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(296,11).Value = "`"10%`""
$ws.Cells.Item(296,12).Value = "Black"
# Row 297
$ws.Cells.Item(297,1).Value = 45599.69183872685
$ws.Cells.Item(297,2).Value = "kusahana8047@gmail.com"
$ws.Cells.Item(297,3).Value = "법학과"
$ws.Cells.Item(297,4).Value = 20192736
$ws.Cells.Item(297,5).Value = "유현우"
$ws.Cells.Item(297,6).Value = "대한민국"
$ws.Cells.Item(297,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(297,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(297,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(297,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(297,11).Value = "`"19.1%`""
$ws.Cells.Item(297,12).Value = "Black"
# Row 298
$ws.Cells.Item(298,1).Value = 45599.699283124995
$ws.Cells.Item(298,2).Value = "ncu11069@naver.com"
$ws.Cells.Item(298,3).Value = "일본학과"
$ws.Cells.Item(298,4).Value = 20231622
$ws.Cells.Item(298,5).Value = "이규민"
$ws.Cells.Item(298,6).Value = "대한민국"
$ws.Cells.Item(298,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(298,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(298,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(298,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(298,11).Value = "`"15%`""
$ws.Cells.Item(298,12).Value = "Red"
# Row 299
$ws.Cells.Item(299,1).Value = 45599.701836875
$ws.Cells.Item(299,2).Value = "wogh1587@naver.com"
$ws.Cells.Item(299,3).Value = "사회복지학과"
$ws.Cells.Item(299,4).Value = 20192366
$ws.Cells.Item(299,5).Value = "현재호"
$ws.Cells.Item(299,6).Value = "대한민국"
$ws.Cells.Item(299,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(299,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(299,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(299,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(299,11).Value = "`"19.1%`""
$ws.Cells.Item(299,12).Value = "Red"
$ws.Cells.Item(299,13).Value = "나는 사후 장기기증에 참여할 뜻이 없다"
# Row 300
$ws.Cells.Item(300,1).Value = 45599.70190223379
$ws.Cells.Item(300,2).Value = "qkhkasin17@naver.com"
$ws.Cells.Item(300,3).Value = "박가현"
$ws.Cells.Item(300,4).Value = 20232953
$ws.Cells.Item(300,5).Value = "박가현"
$ws.Cells.Item(300,6).Value = "대한민국"
$ws.Cells.Item(300,7).Value = "취업자 / 경제활동인구"
$ws.Cells.Item(300,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(300,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(300,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(300,11).Value = "`"19.1%`""
$ws.Cells.Item(300,12).Value = "Red"
# Row 301
$ws.Cells.Item(301,1).Value = 45599.704135277774
$ws.Cells.Item(301,2).Value = "h20231025@glab.hallym.ac.kr"
$ws.Cells.Item(301,3).Value = "국어국문학과"
$ws.Cells.Item(301,4).Value = 20231025
$ws.Cells.Item(301,5).Value = "김지현"
$ws.Cells.Item(301,6).Value = "스페인"
$ws.Cells.Item(301,7).Value = "취업자 / 경제활동인구"
$ws.Cells.Item(301,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(301,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(301,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(301,11).Value = "`"15%`""
$ws.Cells.Item(301,12).Value = "Black"
$ws.Cells.Item(301,14).Value = "나는 사후 장기기증에 참여할 뜻이 있다"
# Row 302
$ws.Cells.Item(302,1).Value = 45599.70781827546
$ws.Cells.Item(302,2).Value = "liz030404@naver.com"
$ws.Cells.Item(302,3).Value = "데이터사이언스학부"
$ws.Cells.Item(302,4).Value = 20243237
$ws.Cells.Item(302,5).Value = "이선주"
$ws.Cells.Item(302,6).Value = "스페인"
$ws.Cells.Item(302,7).Value = "경제활동인구 / 15세이상 인구"
$ws.Cells.Item(302,8).Value = "자기 가구에서 경영하는 농장이나 사업체의 수입을 높이는 데 도운 가족종사자로서 주당 18시간 이상 일한 자"
$ws.Cells.Item(302,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(302,10).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(302,11).Value = "`"19.1%`""
$ws.Cells.Item(302,12).Value = "Red"
$ws.Cells.Item(302,13).Value = "나는 사후 장기기증에 참여할 뜻이 없다"
# Row 303
$ws.Cells.Item(303,1).Value = 45599.71305871528
$ws.Cells.Item(303,2).Value = "8428kyn@naver.com"
$ws.Cells.Item(303,3).Value = "인문학부"
$ws.Cells.Item(303,4).Value = 20241017
$ws.Cells.Item(303,5).Value = "김예나"
$ws.Cells.Item(303,6).Value = "미국"
$ws.Cells.Item(303,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(303,8).Value = "자기 가구에서 경영하는 농장이나 사업체의 수입을 높이는 데 도운 가족종사자로서 주당 18시간 이상 일한 자"
$ws.Cells.Item(303,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(303,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(303,11).Value = "`"15%`""
$ws.Cells.Item(303,12).Value = "Red"
# Row 304
$ws.Cells.Item(304,1).Value = 45599.7187619676
$ws.Cells.Item(304,2).Value = "eojeongmin146@gmail.com"
$ws.Cells.Item(304,3).Value = "법학과"
$ws.Cells.Item(304,4).Value = 20242725
$ws.Cells.Item(304,5).Value = "어정민"
$ws.Cells.Item(304,6).Value = "대한민국"
$ws.Cells.Item(304,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(304,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(304,9).Value = "평균 : 100만원, 중위값 : 200만원"
$ws.Cells.Item(304,10).Value = "평균 : 100만원, 중위값 : 200만원"
$ws.Cells.Item(304,11).Value = "`"10%`""
$ws.Cells.Item(304,12).Value = "Red"
# Row 305
$ws.Cells.Item(305,1).Value = 45599.743392731485
$ws.Cells.Item(305,2).Value = "ann12ann1209@gmail.com"
$ws.Cells.Item(305,3).Value = "경영"
$ws.Cells.Item(305,4).Value = 20222933
$ws.Cells.Item(305,5).Value = "김혜원"
$ws.Cells.Item(305,6).Value = "스페인"
$ws.Cells.Item(305,7).Value = "경제활동인구 / 15세이상 인구"
$ws.Cells.Item(305,8).Value = "조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자"
$ws.Cells.Item(305,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(305,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(305,11).Value = "`"19.1%`""
$ws.Cells.Item(305,12).Value = "Black"
# Row 306
$ws.Cells.Item(306,1).Value = 45599.75039203704
$ws.Cells.Item(306,2).Value = "yxnjxn0331@gmail.com"
$ws.Cells.Item(306,3).Value = "심리학과"
$ws.Cells.Item(306,4).Value = 20217015
$ws.Cells.Item(306,5).Value = "배윤진"
$ws.Cells.Item(306,6).Value = "대한민국"
$ws.Cells.Item(306,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(306,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(306,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(306,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(306,11).Value = "`"19.1%`""
$ws.Cells.Item(306,12).Value = "Red"
$ws.Cells.Item(306,13).Value = "나는 사후 장기기증에 참여할 뜻이 없다"
# Row 307
$ws.Cells.Item(307,1).Value = 45599.751889293984
$ws.Cells.Item(307,2).Value = "alyssa3257@naver.com"
$ws.Cells.Item(307,3).Value = "간호학과"
$ws.Cells.Item(307,4).Value = 20246279
$ws.Cells.Item(307,5).Value = "전영주"
$ws.Cells.Item(307,6).Value = "대한민국"
$ws.Cells.Item(307,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(307,8).Value = "조사대상 주간에 일시적으로 병이 나거나, 날씨가 안 좋거나, 휴가 또는 연가, 노동쟁의 등의 이유로 일을 못한 일시 휴직자"
$ws.Cells.Item(307,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(307,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(307,11).Value = "`"19.1%`""
$ws.Cells.Item(307,12).Value = "Black"
$ws.Cells.Item(307,14).Value = "나는 사후 장기기증에 참여할 뜻이 있다"
# Row 308
$ws.Cells.Item(308,1).Value = 45599.75357353009
$ws.Cells.Item(308,2).Value = "csm06125@naver.com"
$ws.Cells.Item(308,3).Value = "반도체디스플레이"
$ws.Cells.Item(308,4).Value = 20203321
$ws.Cells.Item(308,5).Value = "박근태"
$ws.Cells.Item(308,6).Value = "대한민국"
$ws.Cells.Item(308,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(308,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(308,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(308,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(308,11).Value = "`"15%`""
$ws.Cells.Item(308,12).Value = "Black"
# Row 309
$ws.Cells.Item(309,1).Value = 45599.75710452546
$ws.Cells.Item(309,2).Value = "goemf100@naver.com"
$ws.Cells.Item(309,3).Value = "법학과"
$ws.Cells.Item(309,4).Value = 20222709
$ws.Cells.Item(309,5).Value = "김민정"
$ws.Cells.Item(309,6).Value = "대한민국"
$ws.Cells.Item(309,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(309,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(309,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(309,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(309,11).Value = "`"19.1%`""
$ws.Cells.Item(309,12).Value = "Black"
$ws.Cells.Item(309,14).Value = "나는 사후 장기기증에 참여할 뜻이 있다"
# Row 310
$ws.Cells.Item(310,1).Value = 45599.76175878472
$ws.Cells.Item(310,2).Value = "sin50407899@gmail.com"
$ws.Cells.Item(310,3).Value = "미디어스쿨"
$ws.Cells.Item(310,4).Value = 20202538
$ws.Cells.Item(310,5).Value = "신재화"
$ws.Cells.Item(310,6).Value = "미국"
$ws.Cells.Item(310,7).Value = "경제활동인구 / 15세이상 인구"
$ws.Cells.Item(310,8).Value = "조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자"
$ws.Cells.Item(310,9).Value = "평균 : 100만원, 중위값 : 200만원"
$ws.Cells.Item(310,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(310,11).Value = "`"5%`""
$ws.Cells.Item(310,12).Value = "Red"
# Row 311
$ws.Cells.Item(311,1).Value = 45599.76572497685
$ws.Cells.Item(311,2).Value = "zzun1414@naver.com"
$ws.Cells.Item(311,3).Value = "반도체·디스플레이스쿨"
$ws.Cells.Item(311,4).Value = 20203352
$ws.Cells.Item(311,5).Value = "황준영"
$ws.Cells.Item(311,6).Value = "대한민국"
$ws.Cells.Item(311,7).Value = "취업자 / 15세 이상 인구"
$ws.Cells.Item(311,8).Value = "조사대상 주간에 수입 있는 일을 하지 않았으나, 지난 4주간 구직활동하였으며 즉시 취업 가능한자"
$ws.Cells.Item(311,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(311,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(311,11).Value = "`"19.1%`""
$ws.Cells.Item(311,12).Value = "Red"
# Row 312
$ws.Cells.Item(312,1).Value = 45599.77033122686
$ws.Cells.Item(312,2).Value = "kby5432@naver.com"
$ws.Cells.Item(312,3).Value = "법학과"
$ws.Cells.Item(312,4).Value = 20192737
$ws.Cells.Item(312,5).Value = "윤경빈"
$ws.Cells.Item(312,6).Value = "대한민국"
$ws.Cells.Item(312,7).Value = "경제활동인구 / 15세이상 인구"
$ws.Cells.Item(312,8).Value = "조사 대상 주간 중 수입을 목적으로 1시간 이상 일한 자"
$ws.Cells.Item(312,9).Value = "평균 : 200만원, 중위값 : 100만원"
$ws.Cells.Item(312,10).Value = "평균 : 1,000만원, 중위값 : 100만원"
$ws.Cells.Item(312,11).Value = "`"19.1%`""
$ws.Cells.Item(312,12).Value = "Red"

# ---------------------------------------------------------------------------
# 3) Grow the "Form_Responses1" table to cover the newly added rows.
# ---------------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:N312"))

# ---------------------------------------------------------------------------
# 4) Move the active cell/selection to where it ended up after the edits
#    were made (F316, just past the new data), matching the workbook's
#    last saved cursor position.
# ---------------------------------------------------------------------------
$ws.Range("F316").Select()
